# Updates cryptos list (Coin/Price/Volume columns) to refreshed market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text number format on cells whose new value would otherwise be
# auto-parsed as a number by Excel, so they keep being stored as text
# (matching the source data, e.g. "0.999", "95.21").
$textFmt = "@"

$ws.Range('D2').Value = '43.173.23'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '2.343.15'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('D4').NumberFormat = $textFmt
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = $textFmt
$ws.Range('D5').Value = '302.83'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').NumberFormat = $textFmt
$ws.Range('D6').Value = '95.21'
$ws.Range('E6').Value = '  -1.61%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('E10').Value = '  -2.36%  '
$ws.Range('D11').NumberFormat = $textFmt
$ws.Range('D11').Value = '0.0785'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('D12').NumberFormat = $textFmt
$ws.Range('D12').Value = '18.69'
$ws.Range('E12').Value = '  -3.47%  '
$ws.Range('E13').Value = '  +2.04%  '
$ws.Range('D14').NumberFormat = $textFmt
$ws.Range('D14').Value = '6.77'
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').Value = '2.704.04'
$ws.Range('E15').Value = '  +1.21%  '
$ws.Range('D16').Value = '2.322.20'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').NumberFormat = $textFmt
$ws.Range('D17').Value = '0.801'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('D18').Value = '43.096.17'
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('D19').NumberFormat = $textFmt
$ws.Range('D19').Value = '12.21'
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('E20').Value = '  +3.18%  '
$ws.Range('D21').Value = '0.0₃0891'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = $textFmt
$ws.Range('D22').Value = '68.02'
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('D23').NumberFormat = $textFmt
$ws.Range('D23').Value = '235.79'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = $textFmt
$ws.Range('D24').Value = '2.23'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').NumberFormat = $textFmt
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = $textFmt
$ws.Range('D27').Value = '24.63'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').NumberFormat = $textFmt
$ws.Range('D28').Value = '2.36'
$ws.Range('E28').Value = '  -0.61%  '
$ws.Range('D29').NumberFormat = $textFmt
$ws.Range('D29').Value = '9.18'
$ws.Range('E29').Value = '  +1.17%  '
$ws.Range('D30').NumberFormat = $textFmt
$ws.Range('D30').Value = '31.63'
$ws.Range('E30').Value = '  -2.68%  '
$ws.Range('D31').NumberFormat = $textFmt
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').NumberFormat = $textFmt
$ws.Range('D32').Value = '5.02'
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('D33').NumberFormat = $textFmt
$ws.Range('D33').Value = '0.0727'
$ws.Range('E33').Value = '  +4.23%  '
$ws.Range('D34').NumberFormat = $textFmt
$ws.Range('D34').Value = '17.31'
$ws.Range('E34').Value = '  -2.34%  '
$ws.Range('D35').NumberFormat = $textFmt
$ws.Range('D35').Value = '4.39'
$ws.Range('E35').Value = '  -2.09%  '
$ws.Range('E36').Value = '  +4.05%  '
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = $textFmt
$ws.Range('D39').Value = '22.39'
$ws.Range('E39').Value = '  +20.31%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').NumberFormat = $textFmt
$ws.Range('D40').Value = '2.76'
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('D42').NumberFormat = $textFmt
$ws.Range('D42').Value = '114.03'
$ws.Range('E42').Value = '  -30.40%  '
$ws.Range('D43').Value = '1.937.55'
$ws.Range('E43').Value = '  -2.05%  '
$ws.Range('D44').NumberFormat = $textFmt
$ws.Range('D44').Value = '0.0282'
$ws.Range('E44').Value = '  +1.20%  '
$ws.Range('D45').NumberFormat = $textFmt
$ws.Range('D45').Value = '10.07'
$ws.Range('E45').Value = '  -4.87%  '
$ws.Range('E46').Value = '  +1.72%  '
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').Value = '2.571.14'
$ws.Range('E48').Value = '  +1.27%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = $textFmt
$ws.Range('D49').Value = '53.18'
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = $textFmt
$ws.Range('D50').Value = '2.82'
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('D51').NumberFormat = $textFmt
$ws.Range('D51').Value = '72.15'
$ws.Range('E51').Value = '  +0.08%  '
